$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
